$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 102
$ws.Range("B3").Value = "nm"
$ws.Range("C3").Value = 25

$ws.Range("A3").Select()
